# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-11 from
# 2023-09-15 (45184) to 2023-09-16 (45185), keeping existing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
